# Applies the "Fruta / hortaliza, semanal" update to the Papaya sheet.
#
# Summary of the change (derived from the OOXML diff):
#   - Row 9  (Primera) and Row 10 (Segunda) get refreshed Volumen/Precio/Unidad
#     values (a later weekly observation replaces the one that used to live
#     there), while the rest of their columns (market/product metadata) stay
#     the same.
#   - The observation that used to sit in row 9 is preserved by moving it down
#     into row 11 (which ends up exactly as the former content of row 9).
#   - The observation that used to sit in row 10 is preserved by moving it
#     down into a brand-new row 12 (exactly the former content of row 10).
#   - The observation that used to sit in row 11 is preserved by moving it
#     down into a brand-new row 13 (exactly the former content of row 11).
#
# All of the metadata columns (A, B, C, E, F, G, H, I, J, K, R) are identical
# across rows 9-13, so literal values (taken straight from the known
# before/after cell contents) are written directly instead of relying on
# reading `.Value` back off existing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 9 -- updated in place (Calidad = Primera)
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 44880
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100108
$ws.Cells.Item(9, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(9, 9).Value = 100108004
$ws.Cells.Item(9, 10).Value = "Papaya"
$ws.Cells.Item(9, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 20000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 20000
$ws.Cells.Item(9, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 19).Value = 2000
$ws.Cells.Item(9, 20).Value = 10

# ---------------------------------------------------------------------
# Row 10 -- updated in place (Calidad = Segunda)
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 44880
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100108
$ws.Cells.Item(10, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(10, 9).Value = 100108004
$ws.Cells.Item(10, 10).Value = "Papaya"
$ws.Cells.Item(10, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 12).Value = "Segunda"
$ws.Cells.Item(10, 13).Value = 180
$ws.Cells.Item(10, 14).Value = 15000
$ws.Cells.Item(10, 15).Value = 15000
$ws.Cells.Item(10, 16).Value = 15000
$ws.Cells.Item(10, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 19).Value = 1500
$ws.Cells.Item(10, 20).Value = 10

# ---------------------------------------------------------------------
# Row 11 -- now holds what used to be row 9's observation (Calidad = Primera)
# ---------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 44391
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100108
$ws.Cells.Item(11, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(11, 9).Value = 100108004
$ws.Cells.Item(11, 10).Value = "Papaya"
$ws.Cells.Item(11, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 15
$ws.Cells.Item(11, 14).Value = 1500
$ws.Cells.Item(11, 15).Value = 1500
$ws.Cells.Item(11, 16).Value = 1500
$ws.Cells.Item(11, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(11, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 19).Value = 1500
$ws.Cells.Item(11, 20).Value = 1

# ---------------------------------------------------------------------
# Row 12 (new) -- holds what used to be row 10's observation (Calidad = Segunda)
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(12, 3).Value = "Metropolitana"
$ws.Cells.Item(12, 4).Value = 44391
# Match the "Fecha" date-style formatting used by the rest of column D.
$ws.Cells.Item(12, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 5).Value = 13
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100108
$ws.Cells.Item(12, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(12, 9).Value = 100108004
$ws.Cells.Item(12, 10).Value = "Papaya"
$ws.Cells.Item(12, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 20
$ws.Cells.Item(12, 14).Value = 1000
$ws.Cells.Item(12, 15).Value = 1000
$ws.Cells.Item(12, 16).Value = 1000
$ws.Cells.Item(12, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(12, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 19).Value = 1000
$ws.Cells.Item(12, 20).Value = 1

# ---------------------------------------------------------------------
# Row 13 (new) -- holds what used to be row 11's observation (Calidad = Primera)
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = 9
$ws.Cells.Item(13, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44343
# Match the "Fecha" date-style formatting used by the rest of column D.
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100108
$ws.Cells.Item(13, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(13, 9).Value = 100108004
$ws.Cells.Item(13, 10).Value = "Papaya"
$ws.Cells.Item(13, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 20
$ws.Cells.Item(13, 14).Value = 1700
$ws.Cells.Item(13, 15).Value = 1700
$ws.Cells.Item(13, 16).Value = 1700
$ws.Cells.Item(13, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(13, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(13, 19).Value = 1700
$ws.Cells.Item(13, 20).Value = 1

# Keep selection/active cell consistent with the source workbook.
$ws.Range("A1").Select() | Out-Null
